$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row - copy style from H1 (bold, bordered) onto new header cells
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Data rows
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 3

$ws.Range("I3").Value = 6
$ws.Range("J3").Value = 7

$ws.Range("I4").Value = 4
$ws.Range("J4").Value = 5
